# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The worker detail rows (16-29) get re-sorted: instead of grouping all
# periods for worker 1 (73103546 - JAVIER ROMERO ROMERO GELIS) followed by
# all periods for worker 2 (1143335226 - JAVIER ENRIQUE ROMERO ARRIETA),
# the rows are now interleaved by period (ascending 2405..2411), alternating
# between the two workers, and the "24266" value moves from the first period
# (2411 before the change) to the last period (2411 after the change, i.e.
# rows 28/29) while every other period keeps 52000.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$worker1Doc  = "73103546"
$worker1Name = "JAVIER ROMERO ROMERO GELIS"
$worker2Doc  = "1143335226"
$worker2Name = "JAVIER ENRIQUE ROMERO ARRIETA"

$periods = @("2405","2406","2407","2408","2409","2410","2411")

$row = 16
foreach ($periodo in $periods) {
    $valor = 52000
    if ($periodo -eq "2411") {
        $valor = 24266
    }

    $ws.Cells.Item($row, 3).Value = $worker1Doc
    $ws.Cells.Item($row, 4).Value = $worker1Name
    $ws.Cells.Item($row, 5).Value = $periodo
    $ws.Cells.Item($row, 6).Value = $valor
    $row = $row + 1

    $ws.Cells.Item($row, 3).Value = $worker2Doc
    $ws.Cells.Item($row, 4).Value = $worker2Name
    $ws.Cells.Item($row, 5).Value = $periodo
    $ws.Cells.Item($row, 6).Value = $valor
    $row = $row + 1
}
